$d = $word.ActiveDocument
$p = $d.Paragraphs(1).Range
$p.Collapse(0)
$p.InsertParagraphAfter()
$p.Collapse(0)
$p.Text = "Second line"
